$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-18
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224)
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
